$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column "Price" cells hold numeric-looking text (e.g. "26.124.23",
# "0.5195"); force Text format on each before writing so COM does not
# silently coerce plain decimal-looking values (e.g. "210.09") into
# real numbers - the source data must stay text, matching column D's
# existing cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.124.23"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.667.43"
$ws.Range("E3").Value = "  -1.39%  "
$ws.Range("E4").Value = "  -0.72%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.09"
$ws.Range("E5").Value = "  -4.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5195"
$ws.Range("E6").Value = "  -5.04%  "
$ws.Range("E7").Value = "  -0.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2634"
$ws.Range("E8").Value = "  -4.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06224"
$ws.Range("E9").Value = "  -3.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.15"
$ws.Range("E10").Value = "  -3.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07491"
$ws.Range("E11").Value = "  -2.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.683.45"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.419"
$ws.Range("E13").Value = "  -2.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5592"
$ws.Range("E14").Value = "  -4.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "66.16"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000007899"
$ws.Range("E16").Value = "  -6.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.168.53"
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.782"
$ws.Range("E19").Value = "  -3.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "186.55"
$ws.Range("E20").Value = "  -2.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.35"
$ws.Range("E21").Value = "  -5.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.163"
$ws.Range("E22").Value = "  -1.67%  "
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "147.61"
$ws.Range("E24").Value = "  -1.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1242"
$ws.Range("E25").Value = "  -6.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.557"
$ws.Range("E26").Value = "  -4.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.97"
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06244"
$ws.Range("E28").Value = "  -1.67%  "
$ws.Range("E29").Value = "  -2.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.275"
$ws.Range("E30").Value = "  -4.02%  "
$ws.Range("E31").Value = "  -3.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.419"
$ws.Range("E32").Value = "  -5.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.625"
$ws.Range("E33").Value = "  -3.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9951"
$ws.Range("E34").Value = "  -4.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6028"
$ws.Range("E35").Value = "  -2.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.404"
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("E38").Value = "  -2.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01605"
$ws.Range("E39").Value = "  -1.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.076.27"
$ws.Range("E40").Value = "  -3.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8644"
$ws.Range("E41").Value = "  -1.77%  "
$ws.Range("E42").Value = "  -1.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.16"
$ws.Range("E43").Value = "  -2.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.815.70"
$ws.Range("E44").Value = "  -1.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000110"
$ws.Range("E45").Value = "  +1.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.96"
$ws.Range("E46").Value = "  -2.70%  "
$ws.Range("E47").Value = "  -1.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05252"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.945"
$ws.Range("E49").Value = "  -3.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4249"
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.929"
$ws.Range("E51").Value = "  -4.88%  "
